$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Harvard case classification added: a new "average_doctor" column is inserted,
# shifting the previous "average_doctor" values into "average_doctor_old" (BQ),
# while BP receives newly computed averages. Other avg/variance/stdDev stats were
# recalculated accordingly.

# --- Header row: swap column titles ---
$ws.Range("BP1").Value = "average_doctor_old"
$ws.Range("BQ1").Value = "average_doctor"

# --- Recalculated statistic values (rows 4-13) ---
# Row 4
$ws.Range("E4").Value = 0.488
$ws.Range("F4").Value = 0.052
$ws.Range("G4").Value = 0.229
$ws.Range("N4").Value = 0.465
$ws.Range("O4").Value = 0.057
$ws.Range("P4").Value = 0.239
$ws.Range("Q4").Value = 0.054
$ws.Range("R4").Value = 0.036
$ws.Range("S4").Value = 0.19
$ws.Range("W4").Value = 0.376
$ws.Range("X4").Value = 0.107
$ws.Range("Y4").Value = 0.328
$ws.Range("AI4").Value = 0.404
$ws.Range("AJ4").Value = 0.092
$ws.Range("AK4").Value = 0.304
$ws.Range("AU4").Value = 0.246
$ws.Range("AW4").Value = 0.16
$ws.Range("BA4").Value = 2.018
$ws.Range("BB4").Value = 0.141
$ws.Range("BC4").Value = 0.376
$ws.Range("BG4").Value = 0.715
$ws.Range("BH4").Value = 0.144
$ws.Range("BI4").Value = 0.379
$ws.Range("BM4").Value = 0.744
$ws.Range("BN4").Value = 0.064
$ws.Range("BO4").Value = 0.253
$ws.Range("BP4").Value = 0.673
$ws.Range("BQ4").Value = 0.756
# Row 5
$ws.Range("E5").Value = 0.61
$ws.Range("F5").Value = 0.058
$ws.Range("G5").Value = 0.241
$ws.Range("N5").Value = 0.732
$ws.Range("O5").Value = 0.073
$ws.Range("P5").Value = 0.269
$ws.Range("Q5").Value = 0.035
$ws.Range("R5").Value = 0.015
$ws.Range("S5").Value = 0.122
$ws.Range("W5").Value = 0.339
$ws.Range("X5").Value = 0.099
$ws.Range("Y5").Value = 0.314
$ws.Range("AI5").Value = 0.397
$ws.Range("AJ5").Value = 0.081
$ws.Range("AK5").Value = 0.285
$ws.Range("AU5").Value = 0.461
$ws.Range("AV5").Value = 0.075
$ws.Range("AW5").Value = 0.274
$ws.Range("BA5").Value = 1.308
$ws.Range("BB5").Value = 0.075
$ws.Range("BC5").Value = 0.273
$ws.Range("BG5").Value = 0.382
$ws.Range("BH5").Value = 0.051
$ws.Range("BI5").Value = 0.227
$ws.Range("BM5").Value = 0.528
$ws.Range("BN5").Value = 0.047
$ws.Range("BO5").Value = 0.217
$ws.Range("BP5").Value = 0.436
$ws.Range("BQ5").Value = 0.459
# Row 6
$ws.Range("E6").Value = 0.542
$ws.Range("N6").Value = 0.569
$ws.Range("Q6").Value = 0.042
$ws.Range("W6").Value = 0.357
$ws.Range("AI6").Value = 0.4
$ws.Range("AU6").Value = 0.321
$ws.Range("BA6").Value = 1.581
$ws.Range("BG6").Value = 0.498
$ws.Range("BM6").Value = 0.618
$ws.Range("BP6").Value = 0.527
$ws.Range("BQ6").Value = 0.569
# Row 7
$ws.Range("E7").Value = 0.581
$ws.Range("N7").Value = 0.657
$ws.Range("Q7").Value = 0.038
$ws.Range("W7").Value = 0.346
$ws.Range("AI7").Value = 0.398
$ws.Range("AU7").Value = 0.392
$ws.Range("BA7").Value = 1.404
$ws.Range("BG7").Value = 0.421
$ws.Range("BM7").Value = 0.5610000000000001
$ws.Range("BP7").Value = 0.468
$ws.Range("BQ7").Value = 0.497
# Row 8
$ws.Range("E8").Value = 0.703
$ws.Range("F8").Value = 0.073
$ws.Range("G8").Value = 0.271
$ws.Range("N8").Value = 0.804
$ws.Range("O8").Value = 0.061
$ws.Range("P8").Value = 0.247
$ws.Range("Q8").Value = 0.039
$ws.Range("W8").Value = 0.404
$ws.Range("X8").Value = 0.116
$ws.Range("Y8").Value = 0.341
$ws.Range("AI8").Value = 0.471
$ws.Range("AJ8").Value = 0.132
$ws.Range("AK8").Value = 0.363
$ws.Range("AU8").Value = 0.406
$ws.Range("AV8").Value = 0.083
$ws.Range("AW8").Value = 0.288
$ws.Range("BA8").Value = 1.746
$ws.Range("BB8").Value = 0.11
$ws.Range("BC8").Value = 0.332
$ws.Range("BG8").Value = 0.555
$ws.Range("BH8").Value = 0.112
$ws.Range("BI8").Value = 0.334
$ws.Range("BM8").Value = 0.67
$ws.Range("BN8").Value = 0.063
$ws.Range("BO8").Value = 0.251
$ws.Range("BP8").Value = 0.582
$ws.Range("BQ8").Value = 0.62
# Row 9
$ws.Range("E9").Value = 0.659
$ws.Range("F9").Value = 0.225
$ws.Range("G9").Value = 0.474
$ws.Range("N9").Value = 0.732
$ws.Range("O9").Value = 0.196
$ws.Range("P9").Value = 0.443
$ws.Range("W9").Value = 0.293
$ws.Range("X9").Value = 0.207
$ws.Range("Y9").Value = 0.455
$ws.Range("AI9").Value = 0.415
$ws.Range("AJ9").Value = 0.243
$ws.Range("AK9").Value = 0.493
$ws.Range("BA9").Value = 1.682
$ws.Range("BB9").Value = 0.249
$ws.Range("BC9").Value = 0.499
$ws.Range("BG9").Value = 0.585
$ws.Range("BH9").Value = 0.243
$ws.Range("BI9").Value = 0.493
$ws.Range("BM9").Value = 0.634
$ws.Range("BN9").Value = 0.232
$ws.Range("BO9").Value = 0.482
$ws.Range("BP9").Value = 0.5610000000000001
$ws.Range("BQ9").Value = 0.605
# Row 10
$ws.Range("E10").Value = 0.805
$ws.Range("F10").Value = 0.157
$ws.Range("G10").Value = 0.396
$ws.Range("N10").Value = 0.927
$ws.Range("O10").Value = 0.068
$ws.Range("P10").Value = 0.26
$ws.Range("W10").Value = 0.512
$ws.Range("X10").Value = 0.25
$ws.Range("Y10").Value = 0.5
$ws.Range("AI10").Value = 0.512
$ws.Range("AJ10").Value = 0.25
$ws.Range("AK10").Value = 0.5
$ws.Range("AU10").Value = 0.415
$ws.Range("AV10").Value = 0.243
$ws.Range("AW10").Value = 0.493
$ws.Range("BA10").Value = 2.147
$ws.Range("BB10").Value = 0.217
$ws.Range("BC10").Value = 0.465
$ws.Range("BG10").Value = 0.659
$ws.Range("BH10").Value = 0.225
$ws.Range("BI10").Value = 0.474
$ws.Range("BM10").Value = 0.805
$ws.Range("BN10").Value = 0.157
$ws.Range("BO10").Value = 0.396
$ws.Range("BP10").Value = 0.716
$ws.Range("BQ10").Value = 0.753
# Row 11
$ws.Range("E11").Value = 0.854
$ws.Range("F11").Value = 0.125
$ws.Range("G11").Value = 0.353
$ws.Range("N11").Value = 0.927
$ws.Range("O11").Value = 0.068
$ws.Range("P11").Value = 0.26
$ws.Range("W11").Value = 0.512
$ws.Range("X11").Value = 0.25
$ws.Range("Y11").Value = 0.5
$ws.Range("AI11").Value = 0.585
$ws.Range("AJ11").Value = 0.243
$ws.Range("AK11").Value = 0.493
$ws.Range("AU11").Value = 0.5610000000000001
$ws.Range("AV11").Value = 0.246
$ws.Range("AW11").Value = 0.496
$ws.Range("BA11").Value = 2.147
$ws.Range("BB11").Value = 0.217
$ws.Range("BC11").Value = 0.465
$ws.Range("BG11").Value = 0.659
$ws.Range("BH11").Value = 0.225
$ws.Range("BI11").Value = 0.474
$ws.Range("BM11").Value = 0.805
$ws.Range("BN11").Value = 0.157
$ws.Range("BO11").Value = 0.396
$ws.Range("BP11").Value = 0.716
$ws.Range("BQ11").Value = 0.759
# Row 12
$ws.Range("E12").Value = 1.429
$ws.Range("F12").Value = 0.873
$ws.Range("G12").Value = 0.9350000000000001
$ws.Range("N12").Value = 1.263
$ws.Range("O12").Value = 0.299
$ws.Range("P12").Value = 0.547
$ws.Range("W12").Value = 1.524
$ws.Range("X12").Value = 0.44
$ws.Range("Y12").Value = 0.663
$ws.Range("AI12").Value = 1.625
$ws.Range("AJ12").Value = 1.484
$ws.Range("AK12").Value = 1.218
$ws.Range("AU12").Value = 2.76
$ws.Range("AV12").Value = 3.302
$ws.Range("AW12").Value = 1.817
$ws.Range("BA12").Value = 3.819
$ws.Range("BB12").Value = 0.459
$ws.Range("BC12").Value = 0.678
$ws.Range("BG12").Value = 1.148
$ws.Range("BH12").Value = 0.2
$ws.Range("BI12").Value = 0.448
$ws.Range("BM12").Value = 1.242
$ws.Range("BN12").Value = 0.244
$ws.Range("BO12").Value = 0.494
$ws.Range("BP12").Value = 1.273
$ws.Range("BQ12").Value = 1.255
# Row 13
$ws.Range("E13").Value = 1.394
$ws.Range("F13").Value = 0.266
$ws.Range("G13").Value = 0.516
$ws.Range("N13").Value = 1.78
$ws.Range("O13").Value = 0.431
$ws.Range("P13").Value = 0.656
$ws.Range("W13").Value = 0.978
$ws.Range("X13").Value = 0.196
$ws.Range("Y13").Value = 0.443
$ws.Range("AI13").Value = 1.144
$ws.Range("AJ13").Value = 0.309
$ws.Range("AK13").Value = 0.556
$ws.Range("AU13").Value = 2.059
$ws.Range("AV13").Value = 0.331
$ws.Range("AW13").Value = 0.576
$ws.Range("BA13").Value = 2.21
$ws.Range("BB13").Value = 0.279
$ws.Range("BC13").Value = 0.528
$ws.Range("BG13").Value = 0.55
$ws.Range("BH13").Value = 0.05
$ws.Range("BI13").Value = 0.224
$ws.Range("BM13").Value = 0.796
$ws.Range("BN13").Value = 0.162
$ws.Range("BO13").Value = 0.403
$ws.Range("BP13").Value = 0.737
$ws.Range("BQ13").Value = 0.672
